$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the HMM transition/emission probability inputs.
# All other changed cells in the sheet are formulas derived from these
# four inputs, so they will recalculate automatically.
$ws.Range("B2").Value = 0.8
$ws.Range("C2").Value = 0.2
$ws.Range("H2").Value = 0.9
$ws.Range("I2").Value = 0.1

# Move the selection (and implicitly the scrolled view) to E34, matching
# the author's updated view position after adjusting the results plots.
$ws.Range("E34").Select()
